$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date in column C for rows 2-7 from 2023-09-14 (45183) to 2023-09-15 (45184)
$ws.Range("C2:C7").Value = 45184
